# Daily auto-update: for every data row (2..lastRow), decrement the
# "剩余" (remaining days) counter in column E by 1. When a row's counter
# has reached 1 (i.e. the last day of its cycle), instead of going to 0
# the cycle resets: E is set back to the row's "总天" (total days, column D)
# and the start date in column F advances by that many days.
#
# Row 36 is skipped: its F value ("202510929") is not a valid yyyyMMdd
# date, so it is left untouched (matches source data / diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $totalDays = $ws.Cells.Item($row, 4).Value()   # D: 总天
    $remaining = $ws.Cells.Item($row, 5).Value()   # E: 剩余
    $startDate = $ws.Cells.Item($row, 6).Value()   # F: 开始时间

    if ($remaining -eq $null) {
        continue
    }

    # The start date must be a genuine yyyyMMdd date. Rows with a corrupted
    # date (e.g. row 36's "202510929") are left completely untouched.
    $dateStr = [string][int]$startDate
    $validDate = $true
    try {
        $parsed = [DateTime]::ParseExact($dateStr, "yyyyMMdd", $null)
    } catch {
        $validDate = $false
    }

    if (-not $validDate) {
        continue
    }

    if ($remaining -le 1) {
        $newDate = $parsed.AddDays([int]$totalDays)
        $ws.Cells.Item($row, 6).Value = [int]$newDate.ToString("yyyyMMdd")
        $ws.Cells.Item($row, 5).Value = $totalDays
    } else {
        $ws.Cells.Item($row, 5).Value = $remaining - 1
    }
}
